# AlumniHub.xlsx - "workflow improved for activate email"
#
# The sheet tracked three people (rows 2-4). This change:
#   1. Flips the first person's status from is_student (L2) to
#      is_alumni (K2) -- i.e. the "activate" toggle moves one column left.
#   2. Removes the other two people (Chetan Sonar / Rohan sapkale) and
#      their e-mail hyperlinks entirely, shrinking the sheet back down
#      to a single data row.
#   3. Leaves the first person's own e-mail hyperlink (C2) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-point the boolean flag on row 2 from L2 (is_student) to K2 (is_alumni) ---
$ws.Range("K2").Value = $true
$ws.Range("L2").ClearContents()

# --- 2. Drop the hyperlinks that belong to the rows we are about to remove ---
# (the engine's Hyperlinks.Delete() clears the whole sheet collection, so we
# delete everything and re-create the one hyperlink we need to keep)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:mangesh2003vispute@gmail.com")
$ws.Range("C2").Style = "Hyperlink"

# --- 3. Remove rows 3 and 4 (Chetan Sonar, Rohan sapkale) completely ---
$ws.Rows("3:4").Delete()

# --- 4. Restore the selection to where the author left off ---
[void]$ws.Range("A9").Select()
